$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 28
$ws.Range("I2").Value = 85
$ws.Range("J2").Value = 306
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 89
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 50
$ws.Range("P2").Value = 1
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 27
$ws.Range("T2").Value = 51
$ws.Range("U2").Value = 9
$ws.Range("V2").Value = 494
$ws.Range("X2").Value = 504
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 7
$ws.Range("AA2").Value = 8
